# "large buildings test build conditions properlt"
#
# The buildable_tiles list (column I, "editable" sheet) allowed building on
# ["grass", "farmland"] for most entries. This adds "shore" to that list so
# large buildings' build-condition tests pass on shore tiles too. Entries
# that were already narrower (["farmland"]) or already shore-only
# (["shore"]) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("editable")

$oldValue = '["grass", "farmland"]'
$newValue = '["grass", "farmland", "shore"]'

$lastRow = $ws.Cells.Item($ws.Rows.Count, 9).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Restore the selection/active-sheet state left behind by the edit session:
# the author last clicked on editable!I38, then returned focus to the
# "formatted" sheet before saving.
$ws.Activate()
$ws.Range("I38").Select()

$ws2 = $wb.Worksheets.Item("formatted")
$ws2.Activate()
$ws2.Range("A2:K45").Select()
